# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 06:45"

# --- Update numeric stats for several countries (row = country) ---

# Row 5: India
$ws.Range("B5").Value = 5308014
$ws.Range("C5").Value = 2539
$ws.Range("D5").Value = 4208431
$ws.Range("E5").Value = 1013958

# Row 21: Pakistan
$ws.Range("B21").Value = 305031
$ws.Range("C21").Value = 645
$ws.Range("D21").Value = 292044
$ws.Range("E21").Value = 6572
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 6415

# Row 134: Tailandia
$ws.Range("B134").Value = 3500
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 3338
$ws.Range("E134").Value = 103
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 59

# Row 187: Butan
$ws.Range("B187").Value = 258
$ws.Range("C187").Value = 6
$ws.Range("D187").Value = 186
$ws.Range("E187").Value = 72

# --- Reorder countries: swap "Santa Lucia" (row 204) and "Timor Oriental" (row 205) ---
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# --- Reorder countries: swap "Montserrat" (row 214) and "Islas Malvinas" (row 215), ---
# --- including their associated statistics.                                      ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
